$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values
$ws.Range('D2').Value = '22.456.73'
$ws.Range('D3').Value = '1.572.56'
$ws.Range('D6').Value = "'290.15"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Value = "'49.87"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = "'0.3387"
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Value = "'1.151"
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = "'0.07565"
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').Value = "'21.26"
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = "'6.035"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = "'6.994"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = '1.571.07'
$ws.Range('D17').Value = "'0.00001124"
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Value = "'90.53"
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = "'0.06795"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Value = "'6.363"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D24').Value = '22.457.06'
$ws.Range('D25').Value = "'2.369"
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Value = "'2.673"
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Value = "'20.05"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Value = "'149.53"
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Value = "'5.054"
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Value = "'125.06"
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Value = '1.757.17'
$ws.Range('D32').Value = "'1.066"
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Value = "'6.221"
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Value = "'2.018"
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Value = "'9.839"
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Value = "'0.08389"
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Value = "'0.02481"
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Value = "'0.2304"
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Value = "'1.351"
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Value = "'0.06565"
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Value = "'5.439"
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Value = "'11.33"
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Value = "'0.6256"
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Value = "'14.06"
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Value = "'3.802"
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Value = "'0.5879"
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Value = "'2.075"
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Value = "'127.87"
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').Value = "'0.07304"
$ws.Range('D51').Style = 'Normal'

# Update Volume(1h) column (E) values
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  -0.52%  '
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('E6').Value = '  +0.66%  '
$ws.Range('E7').Value = '  -1.10%  '
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('E9').Value = '  +1.10%  '
$ws.Range('E10').Value = '  +3.58%  '
$ws.Range('E11').Value = '  +2.54%  '
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('E13').Value = '  +3.12%  '
$ws.Range('E14').Value = '  +3.40%  '
$ws.Range('E15').Value = '  +2.33%  '
$ws.Range('E16').Value = '  +0.64%  '
$ws.Range('E17').Value = '  +2.29%  '
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('E21').Value = '  +4.00%  '
$ws.Range('E22').Value = '  +1.72%  '
$ws.Range('E23').Value = '  +3.69%  '
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('E26').Value = '  +4.14%  '
$ws.Range('E27').Value = '  +1.15%  '
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('E31').Value = '  +1.36%  '
$ws.Range('E32').Value = '  +9.44%  '
$ws.Range('E33').Value = '  +6.41%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  +1.60%  '
$ws.Range('E36').Value = '  -0.48%  '
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('E38').Value = '  +2.63%  '
$ws.Range('E39').Value = '  -2.63%  '
$ws.Range('E40').Value = '  +3.52%  '
$ws.Range('E41').Value = '  +2.78%  '
$ws.Range('E42').Value = '  +4.45%  '
$ws.Range('E43').Value = '  +1.53%  '
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('E45').Value = '  +2.07%  '
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('E47').Value = '  +2.83%  '
$ws.Range('E48').Value = '  +2.83%  '
$ws.Range('E49').Value = '  +3.51%  '
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('E51').Value = '  +0.12%  '
